$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2364.375
$ws.Range("I43").Value = 641.0714
$ws.Range("J43").Value = 4777
$ws.Range("K43").Value = 641.0714
$ws.Range("L43").Value = 4777
$ws.Range("M43").Value = -572.0714
$ws.Range("N43").Value = -4915
$ws.Range("H111").Value = 813.1429000000001
$ws.Range("I111").Value = 152
$ws.Range("J111").Value = 2466
$ws.Range("K111").Value = 456
$ws.Range("L111").Value = 7398
$ws.Range("M111").Value = 2611
$ws.Range("N111").Value = -13532
$ws.Range("H112").Value = 1345.7778
$ws.Range("J112").Value = 1345.7778
$ws.Range("L112").Value = 4037.3334
$ws.Range("N112").Value = -6253.3334
$ws.Range("H127").Value = 987.5
$ws.Range("I127").Value = 893.75
$ws.Range("J127").Value = 1175
$ws.Range("K127").Value = 2681.25
$ws.Range("L127").Value = 3525
$ws.Range("M127").Value = 2278.75
$ws.Range("N127").Value = -13445
$ws.Range("H132").Value = 517904.56
$ws.Range("I132").Value = 247037.19
$ws.Range("J132").Value = 911893.4399999999
$ws.Range("K132").Value = 741111.5700000001
$ws.Range("L132").Value = 2735680.32
$ws.Range("M132").Value = -738581.5700000001
$ws.Range("N132").Value = -2740740.32
$ws.Range("H137").Value = 747302.25
$ws.Range("I137").Value = 2168505
$ws.Range("J137").Value = 2862.738
$ws.Range("K137").Value = 6505515
$ws.Range("L137").Value = 8588.214
$ws.Range("M137").Value = -6502965
$ws.Range("N137").Value = -13688.214

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4037.7703
$ws.Range("I32").Value = 4226.5884
$ws.Range("K32").Value = 4226.5884
$ws.Range("M32").Value = -3939.5884
$ws.Range("H132").Value = 3890.077
$ws.Range("I132").Value = 2044.6
$ws.Range("J132").Value = 5043.5
$ws.Range("K132").Value = 6133.799999999999
$ws.Range("L132").Value = 15130.5
$ws.Range("M132").Value = -3603.799999999999
$ws.Range("N132").Value = -20190.5
$ws.Range("H137").Value = 42077.6
$ws.Range("J137").Value = 42077.6
$ws.Range("L137").Value = 42077.6
$ws.Range("N137").Value = -52277.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4446067.5
$ws.Range("I16").Value = 7408777
$ws.Range("J16").Value = 2003.6
$ws.Range("K16").Value = 7408777
$ws.Range("L16").Value = 2003.6
$ws.Range("M16").Value = -7408490
$ws.Range("N16").Value = -2577.6
$ws.Range("H31").Value = 209673.08
$ws.Range("I31").Value = 410314.97
$ws.Range("J31").Value = 2761.125
$ws.Range("K31").Value = 410314.97
$ws.Range("L31").Value = 2761.125
$ws.Range("M31").Value = -410019.97
$ws.Range("N31").Value = -3351.125
$ws.Range("H34").Value = 209673.08
$ws.Range("I34").Value = 410314.97
$ws.Range("J34").Value = 2761.125
$ws.Range("K34").Value = 410314.97
$ws.Range("L34").Value = 2761.125
$ws.Range("M34").Value = -410112.97
$ws.Range("N34").Value = -3165.125
$ws.Range("H113").Value = 4446067.5
$ws.Range("I113").Value = 7408777
$ws.Range("J113").Value = 2003.6
$ws.Range("K113").Value = 7408777
$ws.Range("L113").Value = 2003.6
$ws.Range("M113").Value = -7406607
$ws.Range("N113").Value = -6343.6
$ws.Range("H134").Value = 1725.75
$ws.Range("I134").Value = 1027.1052
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 3081.3156
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -546.3155999999999
$ws.Range("N134").Value = -50070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 79.521736
$ws.Range("I12").Value = 6
$ws.Range("J12").Value = 118.73333
$ws.Range("K12").Value = 18
$ws.Range("L12").Value = 356.19999
$ws.Range("M12").Value = 155
$ws.Range("N12").Value = -702.19999
$ws.Range("H39").Value = 21668
$ws.Range("J39").Value = 21668
$ws.Range("L39").Value = 65004
$ws.Range("N39").Value = -65592
$ws.Range("H109").Value = 3726.6667
$ws.Range("I109").Value = 700
$ws.Range("J109").Value = 4105
$ws.Range("K109").Value = 2100
$ws.Range("L109").Value = 12315
$ws.Range("M109").Value = -1060
$ws.Range("N109").Value = -14395
$ws.Range("H113").Value = 5682670.5
$ws.Range("I113").Value = 671.1
$ws.Range("J113").Value = 10417670
$ws.Range("K113").Value = 2013.3
$ws.Range("L113").Value = 31253010
$ws.Range("M113").Value = 156.6999999999998
$ws.Range("N113").Value = -31257350

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1114.0769
$ws.Range("I97").Value = 1161.2858
$ws.Range("J97").Value = 1059
$ws.Range("K97").Value = 1161.2858
$ws.Range("L97").Value = 1059
$ws.Range("M97").Value = -665.2858000000001
$ws.Range("N97").Value = -2051
$ws.Range("H113").Value = 2537
$ws.Range("I113").Value = 2666.625
$ws.Range("K113").Value = 2666.625
$ws.Range("M113").Value = -496.625
$ws.Range("H132").Value = 3747.7778
$ws.Range("I132").Value = 1982.6666
$ws.Range("J132").Value = 5159.8667
$ws.Range("K132").Value = 5947.9998
$ws.Range("L132").Value = 15479.6001
$ws.Range("M132").Value = -3417.9998
$ws.Range("N132").Value = -20539.6001
$ws.Range("H134").Value = 33307.31
$ws.Range("J134").Value = 33307.31
$ws.Range("L134").Value = 99921.92999999999
$ws.Range("N134").Value = -104991.93
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1138.3572
$ws.Range("I61").Value = 1028.1666
$ws.Range("K61").Value = 1028.1666
$ws.Range("M61").Value = -826.1666
$ws.Range("H113").Value = 1138.3572
$ws.Range("I113").Value = 1028.1666
$ws.Range("K113").Value = 1028.1666
$ws.Range("M113").Value = 1141.8334
$ws.Range("H141").Value = 32188.438
$ws.Range("J141").Value = 32188.438
$ws.Range("L141").Value = 32188.438
$ws.Range("N141").Value = -42548.43799999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 572.875
$ws.Range("I100").Value = 483.2857
$ws.Range("J100").Value = 1200
$ws.Range("K100").Value = 966.5714
$ws.Range("L100").Value = 2400
$ws.Range("M100").Value = -425.5714
$ws.Range("N100").Value = -3482
$ws.Range("H107").Value = 1033.2307
$ws.Range("I107").Value = 626.4
$ws.Range("J107").Value = 1287.5
$ws.Range("K107").Value = 1879.2
$ws.Range("L107").Value = 3862.5
$ws.Range("M107").Value = 40.80000000000018
$ws.Range("N107").Value = -7702.5
$ws.Range("H113").Value = 528.8570999999999
$ws.Range("I113").Value = 520.4
$ws.Range("J113").Value = 550
$ws.Range("K113").Value = 1561.2
$ws.Range("L113").Value = 1650
$ws.Range("M113").Value = 608.8000000000002
$ws.Range("N113").Value = -5990
